$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for rows 2 and 3
$values = @{
    "G"  = -24.69072164948454
    "H"  = -35.56701030927835
    "I"  = -38.83161512027492
    "J"  = -38.83161512027492
    "K"  = -28.4
    "L"  = -48.79725085910653
    "U"  = 75.8
    "V"  = 0.1568384026484585
    "W"  = -1.551912568306011
    "X"  = 0.05709469244900058
    "Y"  = -1.609007260755011
    "Z"  = 0.08672329012069738
    "AA" = -3.367605423930861
    "AB" = 0.05714411778479069
    "AC" = -3.424749541715652
    "AD" = 1.01
    "AE" = 0
    "AF" = 1.01
    "AG" = -74.78999999999999
    "AH" = 0.002085441143069521
    "AI" = 0.01234567901234568
    "AJ" = -0.1830799735624587
    "AK" = -12.44425956738768
    "AM" = -0.113
    "AN" = -0.04508928571428572
    "AP" = 3.338839285714286
    "AQ" = 200
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range("${col}2").Value = $val
    $ws.Range("${col}3").Value = $val
}
